$d = $word.ActiveDocument

# --- Step 1: remove the "Meta description: ..." paragraph that originally sat
#     right under the H1 title (2nd paragraph in the doc). ---
$metaPara = $d.Paragraphs.Item(2)
if ($metaPara.Range.Text -like "*Meta description*") {
    $metaPara.Range.Delete()
}

# --- Step 2: insert a new bold paragraph "Play Genius Free: A Visually Stunning
#     Online Slot" right before the final paragraph (the image-prompt paragraph). ---
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$lastPara.Range.InsertParagraphBefore()

$newPara = $d.Paragraphs.Item($count)
$newPara.Range.Text = "Play Genius Free: A Visually Stunning Online Slot"
$newParaTextRange = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)
$newParaTextRange.Font.Italic = $false
$newParaTextRange.Font.Bold = $true

# --- Step 3: replace the final paragraph's image-prompt text with the new
#     meta-description copy, keeping its existing (italic) run formatting. ---
$old = "Create a feature image that embraces the Arabian Nights-themed gameplay of Genius by Cristaltec. The image should be in a cartoon style and showcase a happy Maya warrior with glasses, standing on a flying Persian carpet. The background should feature a vast desert with the iconic silhouette of an Arabian palace in the distance. The genie from the lamp should be hovering above the Maya warrior, inspiring them to play and interact with the game. The image should be vibrant and colorful, with symbols from the game integrated into the scenery, including the RellWild and 2x symbols. The overall tone of the image should be adventurous and exciting, making potential players curious about the game and portraying the endless opportunities within the reels."
$new = "Read our review of Genius by Cristaltec. Play for free and experience the stunning visuals, Wild and Scatter symbols, and fun gameplay."
$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
